# Auto-generated Word COM-interop script
# Applies text replacements per the target diff using Find/Replace
# across the document's table cells (handles embedded <w:br/> line breaks
# via the '^l' manual-line-break Find/Replace wildcard code).

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Replacement.ClearFormatting()
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
    return $found
}

Replace-Text '- Light yellow oil or brown semi-solid, viscous liquid, or golden yellow solid^l- Odorless resinous oil' '• Solid form with variable appearance: may be a brown amorphous semi-solid, a viscous oil, a chunky golden yellow solid, light yellow oil, brown semi-solid viscous liquid, or an odorless resinous oil.'
Replace-Text '- In water: 2.8 mg/L at 23 °C^l- In 0.15M sodium chloride: 0.77 mg/L at 23 °C^l- Soluble in fixed oils, alcohol, acetone, and glycerol' 'Poor water solubility (~2.8 mg/L at 23 °C and 0.77 mg/L in 0.15 M sodium chloride at 23 °C); soluble in organic solvents with solvent-specific ratios (e.g., “1 part in 1 part” in alcohol and acetone, “1 part in 3 parts” in glycerol); reported as 2.63e-03 g/L in some measurements.'
Replace-Text 'Six different polymorphic forms of dronabinol have been identified. The amorphous form shows approximately 40% more degradation than the polymorphic forms when subjected to stress conditions. The specific crystal systems and thermodynamic data for these polymorphs were not detailed in the available evidence.' 'No validated crystallographic polymorphic forms have been reported for dronabinol. Although literature on related cannabinoids notes polymorphic impurities, no specific polymorph data is available for dronabinol. [Polymorph evidence](https://www.sciencedirect.com/science/article/pii/S0731708524000785)'
Replace-Text 'Stability studies indicate that dronabinol capsules maintain over 97% of their initial Δ9-THC concentration when stored under various conditions (frozen, refrigerated, or at room temperature) for three months. The product packaging and formulation effectively protect against oxidative degradation.' 'Stability studies indicate that impurity levels should be maintained below 1% (optimized conditions achieving even 0.35%); refrigeration (storage between 8 °C and 15 °C) is recommended to counteract instability at room temperature. Detailed assay recovery data by HPLC is not provided. [Stability indicators evidence](https://www.sciencedirect.com/science/article/pii/S0376871611000317)'
Replace-Text 'Specific impurities were not detailed in the provided evidence. However, it is noted that the nitrogen-flushed blister-packaged dronabinol maintained its potency within about 1% of the label claim, while control dosages showed significant degradation.' 'Identified impurities include degradation products such as cannabinol and other unspecified degradants (e.g., labeled as Deg a, Deg b). Additional impurities may arise from excipients, although numerical levels and exact chemical identities are not detailed. [Impurities evidence](https://www.sciencedirect.com/science/article/pii/S0273230024001569)'
Replace-Text 'Dronabinol is classified under BCS Class 2, indicating low solubility and high permeability. The maximum recommended therapeutic daily dose (MRTD) is 0.91 µM/kg/day, and the fraction excreted unchanged in urine is 0.50%.' 'Based on its high lipophilicity (log Kow = 6.97) and extremely low aqueous solubility, dronabinol is categorized as BCS Class II, implying low solubility with high permeability. [Biopharmaceutical classification evidence](https://www.sciencedirect.com/science/article/pii/S0022354923001818)'
Replace-Text ' log Kow = 6.97' ' 6.97 (log Kow)'
Replace-Text 'Experimental findings on moisture absorption were not detailed in the evidence provided.' 'No experimental data on hygroscopic properties or moisture uptake is available from the validated sources.'
Replace-Text 'The chiral properties of dronabinol were not explicitly mentioned in the evidence.' 'The IUPAC name evidences defined stereocenters; however, specific optical rotation values and enantiomeric purity data have not been provided in the available validated data.'
Replace-Text 'Specific degradation temperatures were not provided, but the stability studies suggest that dronabinol is stable at room temperature for extended periods when properly packaged.' 'A specific degradation temperature is not provided. The melting point of 200 °C may be considered an upper thermal limit, but degradation is primarily driven by oxidative factors rather than thermal decomposition.'
Replace-Text 'The glass transition temperature (Tg) was not reported in the available data.' 'No explicit glass transition temperature (Tg) value has been reported; literature mentions its importance in formulation, particularly in relation to freezing processes, but no numerical Tg is available. [Glass transition temperature evidence](https://www.sciencedirect.com/science/article/pii/S0022354924006063)'
Replace-Text '- White to off-white crystalline powder^l- Odorless^l- Bitter taste' '• Appears as a fine crystalline powder that is white to yellowish-white.^l• Odorless and tasteless.'
Replace-Text '- Slightly soluble in water^l- Soluble in acetone^l- Soluble in dimethylformamide^l- Practically insoluble in chloroform and ether' 'No specific solubility data provided.'
Replace-Text '258-259 °C' 'Información no disponible'
Replace-Text 'Acetazolamide exists in two known polymorphic forms, designated as Form A and Form B. Form A is characterized by a monoclinic crystal system, while Form B is triclinic. The thermodynamic stability of these forms indicates that Form B is the more stable polymorph at room temperature, with a transition temperature between 120°C and 148°C. The grinding of Form A can induce a transformation to Form B, which is significant for pharmaceutical formulation processes. [Source](https://www.sciencedirect.com/science/article/pii/S0022286008005115)' 'Acetazolamide exhibits at least two polymorphic forms. The metastable modification I is noted for having a higher density and very high kinetic stability at 20 °C compared to modification II. Both forms can be crystallized from water with only minimal differences in solubility. Although detailed parameters such as exact melting points, crystal systems, or density values are not provided, strong intermolecular hydrogen bonding is identified as the driving force behind these properties. [ScienceDirect](https://www.sciencedirect.com/science/article/pii/S0022354915502724)'
Replace-Text 'Stability studies of acetazolamide oral suspensions have shown that at least 91.2% of the initial concentration remains stable over a 90-day period under various conditions. The pH of the formulations remained stable, and no significant changes in organoleptic properties were observed. [Source](https://pubmed.ncbi.nlm.nih.gov/33214784/)' 'Stability studies using buffered solutions at pH 4 demonstrate that the final dosage forms remain stable for at least 90 days at 37 °C with a potency loss of only 5%. Additionally, FDA guidelines support a tentative expiry of 2 years at 25 °C. Stability-indicating HPLC methods confirm robust separation (resolution >2) between acetazolamide and its degradation products, with a mass balance close to 99.6%. [ScienceDirect](https://www.sciencedirect.com/science/article/pii/S0731708509007377)'
Replace-Text 'The stability-indicating LC method has identified several process-related impurities, including imp-1, imp-2, imp-3, and imp-4, with purities exceeding 99%. The method has demonstrated a mass balance close to 99.6%, indicating effective separation and quantification of acetazolamide and its impurities. [Source](https://www.sciencedirect.com/science/article/pii/S0731708509007377)' 'Degradation and stability studies have identified process-related impurities. The reference standard of acetazolamide is reported at 99.1% purity, while related impurities (imp-1, imp-2, imp-3, and imp-4) range from 99.4% to 99.7% purity. Specific CAS numbers and chemical structures for these impurities were not provided. [ScienceDirect](https://www.sciencedirect.com/science/article/pii/S0731708509007377)'
Replace-Text 'Acetazolamide''s solubility and permeability characteristics have been reviewed in the context of the Biopharmaceutical Classification System (BCS). The available data suggest that acetazolamide does not meet the criteria for a biowaiver due to insufficient evidence regarding its solubility and permeability. [Source](https://www.sciencedirect.com/science/article/pii/S0022354916326922)' 'Based on the Biopharmaceutics Classification System (BCS), acetazolamide cannot be definitively classified due to insufficient solubility and permeability data. The lack of conclusive in vitro and in vivo absorption studies necessitates bioequivalence testing and precludes a biowaiver. [ScienceDirect](https://www.sciencedirect.com/science/article/pii/S0022354916326922)'
Replace-Text ' -0.3' ' –0.45'
Replace-Text 'Information regarding the hygroscopic nature of acetazolamide is not available, which is critical for understanding its stability in various formulations.' 'No experimental data on hygroscopicity or moisture uptake are available. Further research is required to elucidate its impact on formulation and storage stability.'
Replace-Text 'There is no available data on the chiral properties or specific optical rotation of acetazolamide.' 'No data on chirality or specific optical rotation have been reported. Additional stereochemical studies are warranted.'
Replace-Text 'The degradation of acetazolamide occurs significantly under acidic and basic conditions, with specific degradation temperatures not explicitly defined in the available literature.' 'While degradation under hydrolytic stress conditions has been observed, specific degradation temperature thresholds or kinetic data are not provided. Additional research is needed to determine these parameters.'
Replace-Text 'Information regarding the glass transition temperature (Tg) of acetazolamide is not available.' 'No data regarding the glass transition temperature (Tg) determined by techniques such as DSC are available. Further studies are recommended.'
Replace-Text ' Not applicable (decomposes)' ' Información no disponible'

Write-Host "Done."
